$d = $word.ActiveDocument

$replacements = @(
  "What are readings and measurements and the absolute uncertainty of both? And to what decimal place should they be?",
  "They should be to no greater no. of decimal places than the value it is for. E.g., if you calculate the uncertainty for 0.29mm to be 0.0081mm then it must be 0.01mm.",
  "Look at the lowest number of significant figures (e.g., a reading of e.g., 2.3 x 10",
  "Choose the furthest gradient (yet, in most cases, they may appear symmetrical) and the difference between the best and worst is the uncertainty divided by the line of best fit ",
  "What should the uncertainty be taken as when calculating the mean from a set of data?",
  "E.g., an apple of 1N (0.1kg) would have an order of magnitude of 10"
)

foreach ($text in $replacements) {
    $d.Content.Find.Execute($text, $true, $false, $false, $false, $false,
                             $true, 1, $false, $text, 2)
}
